$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows(1).Font.Bold = $true
$ws1.Rows(1).RowHeight = $ws1.Rows(1).RowHeight
Write-Host "done"
